$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = 'The old chapel relied on donations from the local community. The Alder family had donated a great deal of money to the chapel and always considered themselves better because of it. However  it still wasn’t enough to match the amount of money that the Dwight family had managed to raise for the chapel. This rivalry had been going on for years and had started to fray relations between everyone in the village. However, a new property developer that had just moved in wanted to grease the wheels of his business plans and so started to spend generously in the village. When Mr. Roberts donated to the chapel, he made sure he at least donated more than the Alder Family.'
$ws.Range("E9").Value = 'The old chapel doors stood wide open allowing full view of its interior. Inside there stood a magnificent alter that had been in the church for hundreds of years. To the left of the alter stood the pulpit where the pastor would often hold his sermons from. A recent Archaeological discovery has lead researchers to believe that there is a secret underground chamber somewhere to the left of the alter in the church. The local community considered the church an integral part of its identity, though the number of people attending had dwindled significantly over the past few years. The old chapel relied on donations from the local community. The Alder family had donated a great deal of money to the chapel and always considered themselves better because of it. However  it still wasn’t enough to match the amount of money that the Dwight family had managed to raise for the chapel. This rivalry had been going on for years and had started to fray relations between everyone in the village. However, a new property developer that had just moved in wanted to grease the wheels of his business plans and so started to spend generously in the village. When Mr. Roberts donated to the chapel, he made sure he at least donated more than the Alder Family.'
$ws.Range("J16").Value = 'C.There is not enough information'
$ws.Range("J17").Value = 'B. Down the mountain'
$ws.Range("J24").Value = 'C.There is not enough information'
$ws.Range("J25").Value = 'B.Back towards the surgery'
$ws.Range("J32").Value = 'C.There is not enough information'
$ws.Range("J33").Value = 'B.Right'
$ws.Range("J40").Value = 'C.There is not enough information'
$ws.Range("J41").Value = 'A.Up'
$ws.Range("J48").Value = 'C.There is not enough information'
$ws.Range("J49").Value = 'A.Towards the bridge'
$ws.Range("J56").Value = 'C.There is not enough information'
$ws.Range("J57").Value = 'B.Right'
$ws.Range("J64").Value = 'C.There is not enough information'
$ws.Range("J65").Value = 'A.The wildlife charity '
$ws.Range("J72").Value = 'C.There is not enough information'
$ws.Range("J73").Value = 'A.Up '
$ws.Range("J80").Value = 'C.There is not enough information'
$ws.Range("J81").Value = 'B.Left '
$ws.Range("J88").Value = 'C.There is not enough information'
$ws.Range("J89").Value = 'B. Back '
$ws.Range("J96").Value = 'C.There is not enough information'
$ws.Range("J97").Value = 'A.The Helicopter'
$ws.Range("J104").Value = 'C.There is not enough information'
$ws.Range("J105").Value = 'A.To the Left '
$ws.Range("J112").Value = 'C.There is not enough information'
$ws.Range("J113").Value = 'B.In front of the lathe.'
$ws.Range("J120").Value = 'C.There is not enough information'
$ws.Range("J121").Value = 'B.Down '
$ws.Range("J128").Value = 'C.There is not enough information'
$ws.Range("J129").Value = 'B.Right'
$ws.Range("J136").Value = 'C.There is not enough information'
$ws.Range("J137").Value = 'A.Team Blue'
$ws.Range("J144").Value = 'C.There is not enough information'
$ws.Range("J145").Value = 'B.Down the hill '
$ws.Range("J152").Value = 'C.There is not enough information'
$ws.Range("J153").Value = 'A.Left'
$ws.Range("J160").Value = 'C.There is not enough information'
$ws.Range("J161").Value = 'A.The King Fisher'
$ws.Range("J168").Value = 'C.There is not enough information'
$ws.Range("J169").Value = 'B.Down'
$ws.Range("J176").Value = 'C.There is not enough information'
$ws.Range("J177").Value = 'A.Left'
$ws.Range("J184").Value = 'C.There is not enough information'
$ws.Range("J185").Value = 'A.Towards the Sculpture of the Whale'
$ws.Range("J192").Value = 'C.There is not enough information'
$ws.Range("J193").Value = 'A.Oak'
